# Update Name of Algo
# Apply updated numeric values to the RandomForest imputation result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -7.366099999999996
$ws.Range("B7").Value = 4.578199999999999
$ws.Range("A8").Value = -22.36850000000002
$ws.Range("A10").Value = -21.78439999999999
$ws.Range("A12").Value = -21.59330000000002
$ws.Range("B15").Value = 4.566099999999995
$ws.Range("A18").Value = -22.4215
$ws.Range("B18").Value = 4.305199999999997
$ws.Range("D18").Value = -8.249799999999988
$ws.Range("D19").Value = -8.589899999999991
$ws.Range("B20").Value = 9.182399999999998
$ws.Range("D27").Value = -8.791999999999998
$ws.Range("B29").Value = 4.875800000000002
$ws.Range("B30").Value = 5.465500000000002
$ws.Range("B31").Value = 5.555800000000005
$ws.Range("D31").Value = -8.187299999999999
$ws.Range("A37").Value = -20.4952
$ws.Range("D38").Value = -8.462099999999994
$ws.Range("B40").Value = 9.033599999999996
$ws.Range("D42").Value = -8.752399999999996
$ws.Range("D44").Value = -7.452899999999999
$ws.Range("D47").Value = -7.516800000000003
$ws.Range("B50").Value = 4.718
$ws.Range("A55").Value = -21.6213
$ws.Range("D58").Value = -8.432699999999993
$ws.Range("D65").Value = -7.440899999999996
$ws.Range("A68").Value = -21.47660000000001
$ws.Range("B68").Value = 4.528899999999997
$ws.Range("D73").Value = -7.668299999999995
$ws.Range("B76").Value = 6.3453
$ws.Range("A77").Value = -20.81079999999999
$ws.Range("A78").Value = -20.02129999999998
$ws.Range("A81").Value = -22.00650000000001
$ws.Range("A82").Value = -21.6487
$ws.Range("B87").Value = 4.594099999999994
$ws.Range("B88").Value = 4.543599999999997
$ws.Range("D90").Value = -7.876900000000003
$ws.Range("D94").Value = -6.945399999999993
$ws.Range("D95").Value = -8.1412
$ws.Range("B96").Value = 5.319100000000005
$ws.Range("B98").Value = 5.874499999999999
$ws.Range("B101").Value = 8.832299999999998
$ws.Range("D101").Value = -7.703599999999999
$ws.Range("B102").Value = 8.675400000000005
